$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 10).Value = 4.6
$ws.Cells.Item(2, 13).Value = 1.33
$ws.Cells.Item(2, 14).Value = 3.05
$ws.Cells.Item(2, 15).Value = 1.88
$ws.Cells.Item(2, 16).Value = 1.72
$ws.Cells.Item(2, 17).Value = 3.45
$ws.Cells.Item(2, 18).Value = 1.26
$ws.Cells.Item(2, 21).Value = 1.96
$ws.Cells.Item(2, 22).Value = 1.75
$ws.Cells.Item(2, 31).Value = 13.5
$ws.Cells.Item(2, 34).Value = 5.6
$ws.Cells.Item(2, 35).Value = 6.6
$ws.Cells.Item(2, 37).Value = 10.75
$ws.Cells.Item(2, 38).Value = 11.75
$ws.Cells.Item(3, 13).Value = 1.31
$ws.Cells.Item(3, 14).Value = 3.15
$ws.Cells.Item(3, 17).Value = 3.35
$ws.Cells.Item(3, 18).Value = 1.28
$ws.Cells.Item(3, 21).Value = 1.78
$ws.Cells.Item(3, 22).Value = 1.93
$ws.Cells.Item(4, 13).Value = 1.35
$ws.Cells.Item(4, 14).Value = 2.95
$ws.Cells.Item(4, 17).Value = 3.65
$ws.Cells.Item(4, 18).Value = 1.24
$ws.Cells.Item(4, 21).Value = 1.85
$ws.Cells.Item(4, 22).Value = 1.85
$ws.Cells.Item(5, 7).Value = 3.2
$ws.Cells.Item(5, 9).Value = 2.15
$ws.Cells.Item(5, 10).Value = 3.6
$ws.Cells.Item(5, 11).Value = 2.25
$ws.Cells.Item(5, 12).Value = 2.75
$ws.Cells.Item(5, 13).Value = 1.22
$ws.Cells.Item(5, 14).Value = 4
$ws.Cells.Item(5, 15).Value = 1.73
$ws.Cells.Item(5, 16).Value = 2.08
$ws.Cells.Item(5, 17).Value = 2.75
$ws.Cells.Item(5, 18).Value = 1.4
$ws.Cells.Item(5, 23).Value = 12
$ws.Cells.Item(5, 24).Value = 17
$ws.Cells.Item(5, 37).Value = 21
$ws.Cells.Item(5, 38).Value = 17
$ws.Cells.Item(6, 7).Value = 2.25
$ws.Cells.Item(6, 8).Value = 3.25
$ws.Cells.Item(6, 10).Value = 2.85
$ws.Cells.Item(6, 11).Value = 2.05
$ws.Cells.Item(6, 12).Value = 3.55
$ws.Cells.Item(6, 13).Value = 1.27
$ws.Cells.Item(6, 14).Value = 3.1
$ws.Cells.Item(6, 15).Value = 1.8
$ws.Cells.Item(6, 16).Value = 1.8
$ws.Cells.Item(6, 17).Value = 2.85
$ws.Cells.Item(6, 19).Value = 1.39
$ws.Cells.Item(6, 21).Value = 1.65
$ws.Cells.Item(6, 22).Value = 2
$ws.Cells.Item(6, 23).Value = 8.25
$ws.Cells.Item(6, 24).Value = 11.5
$ws.Cells.Item(6, 25).Value = 8.75
$ws.Cells.Item(6, 26).Value = 23
$ws.Cells.Item(6, 27).Value = 17.5
$ws.Cells.Item(6, 28).Value = 26
$ws.Cells.Item(6, 29).Value = 10
$ws.Cells.Item(6, 30).Value = 6.3
$ws.Cells.Item(6, 31).Value = 13
$ws.Cells.Item(6, 32).Value = 55
$ws.Cells.Item(6, 33).Value = 400
$ws.Cells.Item(6, 34).Value = 9.75
$ws.Cells.Item(6, 35).Value = 16
$ws.Cells.Item(6, 36).Value = 10.5
$ws.Cells.Item(6, 37).Value = 37
$ws.Cells.Item(6, 38).Value = 25
$ws.Cells.Item(6, 39).Value = 32
$ws.Cells.Item(7, 7).Value = 3.2
$ws.Cells.Item(7, 8).Value = 3.45
$ws.Cells.Item(7, 9).Value = 2.05
$ws.Cells.Item(7, 10).Value = 3.65
$ws.Cells.Item(7, 12).Value = 2.65
$ws.Cells.Item(7, 23).Value = 10.75
$ws.Cells.Item(7, 24).Value = 17.5
$ws.Cells.Item(7, 25).Value = 11.25
$ws.Cells.Item(7, 27).Value = 27
$ws.Cells.Item(7, 32).Value = 60
$ws.Cells.Item(7, 34).Value = 7.9
$ws.Cells.Item(7, 35).Value = 10.25
$ws.Cells.Item(7, 37).Value = 19
$ws.Cells.Item(7, 38).Value = 16
$ws.Cells.Item(8, 7).Value = 2
$ws.Cells.Item(8, 8).Value = 3.35
$ws.Cells.Item(8, 9).Value = 3.2
$ws.Cells.Item(8, 10).Value = 2.6
$ws.Cells.Item(8, 12).Value = 3.7
$ws.Cells.Item(8, 23).Value = 6.4
$ws.Cells.Item(8, 24).Value = 8.25
$ws.Cells.Item(8, 25).Value = 7.3
$ws.Cells.Item(8, 26).Value = 14.5
$ws.Cells.Item(8, 32).Value = 45
$ws.Cells.Item(8, 34).Value = 8.5
$ws.Cells.Item(8, 35).Value = 14
$ws.Cells.Item(8, 37).Value = 32
$ws.Cells.Item(9, 21).Value = 1.98
$ws.Cells.Item(9, 22).Value = 1.74
$ws.Cells.Item(11, 15).Value = 1.65
$ws.Cells.Item(11, 16).Value = 2.2
$ws.Cells.Item(11, 17).Value = 2.5
$ws.Cells.Item(11, 18).Value = 1.5
$ws.Cells.Item(12, 7).Value = 2.3
$ws.Cells.Item(12, 9).Value = 2.75
$ws.Cells.Item(12, 12).Value = 3.2
$ws.Cells.Item(12, 33).Value = 126
$ws.Cells.Item(12, 34).Value = 13
$ws.Cells.Item(12, 35).Value = 17
$ws.Cells.Item(12, 39).Value = 26
$ws.Cells.Item(13, 7).Value = 2.35
$ws.Cells.Item(13, 8).Value = 3.2
$ws.Cells.Item(13, 9).Value = 2.8
$ws.Cells.Item(13, 10).Value = 3
$ws.Cells.Item(13, 13).Value = 1.29
$ws.Cells.Item(13, 14).Value = 3.5
$ws.Cells.Item(13, 15).Value = 1.93
$ws.Cells.Item(13, 16).Value = 1.88
$ws.Cells.Item(13, 17).Value = 3.25
$ws.Cells.Item(13, 18).Value = 1.33
$ws.Cells.Item(13, 24).Value = 12
$ws.Cells.Item(13, 25).Value = 10
$ws.Cells.Item(13, 26).Value = 23
$ws.Cells.Item(13, 31).Value = 13
$ws.Cells.Item(13, 33).Value = 500
$ws.Cells.Item(13, 40).Value = 1.05
$ws.Cells.Item(13, 41).Value = 8.5
$ws.Cells.Item(14, 7).Value = 1.65
$ws.Cells.Item(14, 9).Value = 5.75
$ws.Cells.Item(14, 10).Value = 2.38
$ws.Cells.Item(14, 12).Value = 6.5
$ws.Cells.Item(14, 15).Value = 2.4
$ws.Cells.Item(14, 16).Value = 1.53
$ws.Cells.Item(14, 19).Value = 1.53
$ws.Cells.Item(14, 20).Value = 2.38
$ws.Cells.Item(14, 21).Value = 2.38
$ws.Cells.Item(14, 22).Value = 1.53
$ws.Cells.Item(14, 23).Value = 5
$ws.Cells.Item(14, 24).Value = 6.5
$ws.Cells.Item(14, 26).Value = 12
$ws.Cells.Item(14, 28).Value = 41
$ws.Cells.Item(14, 29).Value = 7
$ws.Cells.Item(14, 31).Value = 23
$ws.Cells.Item(14, 35).Value = 29
$ws.Cells.Item(14, 36).Value = 21
$ws.Cells.Item(14, 37).Value = 67
$ws.Cells.Item(14, 39).Value = 67
$ws.Cells.Item(14, 40).Value = 1.1
$ws.Cells.Item(14, 41).Value = 7
$ws.Cells.Item(14, 44).Value = 1.83
$ws.Cells.Item(14, 45).Value = 1.98
$ws.Cells.Item(16, 7).Value = 1.8
$ws.Cells.Item(16, 8).Value = 3.7
$ws.Cells.Item(16, 9).Value = 4.2
$ws.Cells.Item(16, 10).Value = 2.4
$ws.Cells.Item(16, 12).Value = 4.75
$ws.Cells.Item(16, 15).Value = 1.85
$ws.Cells.Item(16, 16).Value = 1.95
$ws.Cells.Item(16, 19).Value = 1.36
$ws.Cells.Item(16, 20).Value = 3
$ws.Cells.Item(16, 24).Value = 8.5
$ws.Cells.Item(16, 33).Value = 251
$ws.Cells.Item(16, 34).Value = 13
$ws.Cells.Item(16, 36).Value = 15
$ws.Cells.Item(17, 40).Value = 1.07
$ws.Cells.Item(17, 41).Value = 9
$ws.Cells.Item(18, 7).Value = 1.34
$ws.Cells.Item(18, 8).Value = 4.3
$ws.Cells.Item(18, 9).Value = 9
$ws.Cells.Item(18, 10).Value = 1.87
$ws.Cells.Item(18, 11).Value = 2.18
$ws.Cells.Item(18, 12).Value = 8
$ws.Cells.Item(18, 13).Value = 1.32
$ws.Cells.Item(18, 14).Value = 2.85
$ws.Cells.Item(18, 15).Value = 1.93
$ws.Cells.Item(18, 16).Value = 1.7
$ws.Cells.Item(18, 17).Value = 3.15
$ws.Cells.Item(18, 18).Value = 1.26
$ws.Cells.Item(18, 19).Value = 1.42
$ws.Cells.Item(18, 20).Value = 2.47
$ws.Cells.Item(18, 21).Value = 2.37
$ws.Cells.Item(18, 22).Value = 1.45
$ws.Cells.Item(18, 23).Value = 5.2
$ws.Cells.Item(18, 26).Value = 7.7
$ws.Cells.Item(18, 27).Value = 13.5
$ws.Cells.Item(18, 29).Value = 8.75
$ws.Cells.Item(18, 31).Value = 29
$ws.Cells.Item(18, 34).Value = 18
$ws.Cells.Item(18, 35).Value = 60
$ws.Cells.Item(18, 37).Value = 300
$ws.Cells.Item(19, 41).Value = 5.8
$ws.Cells.Item(20, 7).Value = 6.7
$ws.Cells.Item(20, 8).Value = 3.9
$ws.Cells.Item(20, 9).Value = 1.47
$ws.Cells.Item(20, 10).Value = 6.4
$ws.Cells.Item(20, 11).Value = 2.15
$ws.Cells.Item(20, 12).Value = 2.02
$ws.Cells.Item(20, 13).Value = 1.31
$ws.Cells.Item(20, 14).Value = 2.87
$ws.Cells.Item(20, 15).Value = 1.93
$ws.Cells.Item(20, 16).Value = 1.7
$ws.Cells.Item(20, 17).Value = 3.15
$ws.Cells.Item(20, 18).Value = 1.26
$ws.Cells.Item(20, 19).Value = 1.42
$ws.Cells.Item(20, 20).Value = 2.47
$ws.Cells.Item(20, 21).Value = 2.1
$ws.Cells.Item(20, 22).Value = 1.57
$ws.Cells.Item(20, 23).Value = 15
$ws.Cells.Item(20, 24).Value = 40
$ws.Cells.Item(20, 25).Value = 22
$ws.Cells.Item(20, 26).Value = 150
$ws.Cells.Item(20, 27).Value = 90
$ws.Cells.Item(20, 28).Value = 90
$ws.Cells.Item(20, 29).Value = 9
$ws.Cells.Item(20, 30).Value = 7.8
$ws.Cells.Item(20, 31).Value = 22
$ws.Cells.Item(20, 34).Value = 5.6
$ws.Cells.Item(20, 35).Value = 6
$ws.Cells.Item(20, 36).Value = 8.5
$ws.Cells.Item(20, 37).Value = 9.5
$ws.Cells.Item(20, 38).Value = 13.5
$ws.Cells.Item(20, 39).Value = 35
$ws.Cells.Item(21, 7).Value = 1.31
$ws.Cells.Item(21, 8).Value = 4.55
$ws.Cells.Item(21, 10).Value = 1.8
$ws.Cells.Item(21, 11).Value = 2.25
$ws.Cells.Item(21, 12).Value = 8.25
$ws.Cells.Item(21, 13).Value = 1.28
$ws.Cells.Item(21, 14).Value = 3
$ws.Cells.Item(21, 16).Value = 1.78
$ws.Cells.Item(21, 17).Value = 2.95
$ws.Cells.Item(21, 18).Value = 1.3
$ws.Cells.Item(21, 19).Value = 1.39
$ws.Cells.Item(21, 20).Value = 2.57
$ws.Cells.Item(21, 21).Value = 2.32
$ws.Cells.Item(21, 22).Value = 1.47
$ws.Cells.Item(21, 23).Value = 5.5
$ws.Cells.Item(21, 24).Value = 5.3
$ws.Cells.Item(21, 26).Value = 7.4
$ws.Cells.Item(21, 27).Value = 12.5
$ws.Cells.Item(21, 28).Value = 40
$ws.Cells.Item(21, 29).Value = 9.5
$ws.Cells.Item(21, 30).Value = 9.5
$ws.Cells.Item(21, 31).Value = 30
$ws.Cells.Item(21, 32).Value = 200
$ws.Cells.Item(21, 34).Value = 19.5
$ws.Cells.Item(21, 37).Value = 300
$ws.Cells.Item(21, 38).Value = 150
$ws.Cells.Item(21, 39).Value = 150
